# LH_TC_USERHOME_REVIEWS.xlsx — "Add files via upload / v1.1 Close owner status"
#
# Hala Eldaly reviewed + closed out the three open review items on the
# "LH-TC-REGISTERATION-Reviews" sheet (Owner Status: Open -> Closed) and
# logged the closure as a new v1.1 entry on the "Version History" sheet.
# Also fixes the capitalisation of the v1.0 "Updated Section" note and
# leaves the workbook focused/selected the way it was left in the author's
# last save.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LH-TC-REGISTERATION-Reviews")
$ws2 = $wb.Worksheets.Item("Version History")

# --- 1. Close out the three review rows on the Reviews sheet ----------------
$ws1.Range("I2").Value = "Closed"
$ws1.Range("I3").Value = "Closed"
$ws1.Range("I4").Value = "Closed"

# Give the "Reviewer verification" cell for the last row (J4) the same
# centred box styling the rest of that column/row uses, but top-aligned
# (matches the new cellXfs entry added for this edit).
$srcFmt = $ws1.Range("I2")
$dstFmt = $ws1.Range("J4")
[void]$srcFmt.Copy()
[void]$dstFmt.PasteSpecial(-4122)
$dstFmt.VerticalAlignment = -4160

# --- 2. Fix capitalisation of the existing v1.0 Version History note -------
$ws2.Range("C2").Value = "Review  the user home feature test cases"

# --- 3. Log the new v1.1 "Close owner status" version-history entry --------
$fmtSrc = $ws2.Range("A2:D2")
$fmtDst = $ws2.Range("A3:D3")
[void]$fmtSrc.Copy()
[void]$fmtDst.PasteSpecial(-4122)

$ws2.Range("A3").Value = "v1.1"
$ws2.Range("B3").Value = "Hala Eldaly"
$ws2.Range("C3").Value = "Close owner status"
$ws2.Range("D3").Value = "13/5/2025"

# --- 4. Restore the selection/active-sheet state from the last save --------
[void]$ws2.Range("C4").Select()
[void]$ws1.Activate()
[void]$ws1.Range("G8").Select()
